$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the Intervention=1 flag that was mistakenly entered for the
# St.Gallen control row (row 49).
$ws.Range("C49").ClearContents() | Out-Null

# New rows of recruitment data for the newly added centers/visits.
# Copy the date-formatted style from an existing date cell (A52) so the
# new date cells reuse the same cell style instead of creating new ones.
$xlPasteFormats = -4122

# Row 53: Zuerich, Control
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A53").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A53").Value = 45818
$ws.Range("B53").Value = "Zuerich"
$ws.Range("H53").Value = 1

# Row 54: St. Gallen, Intervention / Uptake / E-cigarettes
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A54").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A54").Value = 45819
$ws.Range("B54").Value = "St. Gallen"
$ws.Range("C54").Value = 1
$ws.Range("D54").Value = 1
$ws.Range("E54").Value = 1

# Row 55: Zuerich, Intervention / Uptake (declined)
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A55").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A55").Value = 45824
$ws.Range("B55").Value = "Zuerich"
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 0

# Row 56: Zuerich, Intervention / Uptake (declined)
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A56").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A56").Value = 45824
$ws.Range("B56").Value = "Zuerich"
$ws.Range("C56").Value = 1
$ws.Range("D56").Value = 0

$excel.CutCopyMode = $false

# Restore the view roughly where the editor last left it (scrolled down,
# selection parked past the last row).
$ws.Range("A14").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("I58").Select() | Out-Null
